# Updates crypto price/volume figures (and re-orders the Frax / BabyDogeCoin
# rows) to match the latest scrape, as produced by the scheduled GitHub
# Actions job that refreshes cryptos.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: cells in the "Price" column hold numeric-looking text (e.g. "218.63").
# Excel's COM layer auto-converts such strings to real numbers on assignment,
# which would both change the cell type and silently drop things like
# trailing zeros (e.g. "0.8660" -> 0.866). Forcing NumberFormat to Text ("@")
# immediately before writing those values keeps them as text, matching the
# original inline-string cells.

$ws.Range('D2').Value = '26.282.64'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').Value = '1.664.07'
$ws.Range('E3').Value = '  +0.57%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.63'
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5324'
$ws.Range('E6').Value = '  +1.37%  '
$ws.Range('E7').Value = '  +0.79%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2645'
$ws.Range('E8').Value = '  +1.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06375'
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('E10').Value = '  +0.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07815'
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.565'
$ws.Range('E12').Value = '  +1.22%  '
$ws.Range('D13').Value = '1.665.10'
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('D14').Value = '1.892.67'
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5522'
$ws.Range('E15').Value = '  +0.92%  '
$ws.Range('D16').Value = '0.0₅8209'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.62'
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.689'
$ws.Range('E19').Value = '  +2.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.41'
$ws.Range('E20').Value = '  +1.07%  '
$ws.Range('E21').Value = '  +1.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.029'
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.59'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1230'
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.194'
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.12'
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('E28').Value = '  +3.65%  '
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.281'
$ws.Range('E30').Value = '  +0.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.606'
$ws.Range('E31').Value = '  +2.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.278'
$ws.Range('E32').Value = '  +0.67%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9618'
$ws.Range('E34').Value = '  +1.18%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.825'
$ws.Range('E35').Value = '  +1.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.417'
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5804'
$ws.Range('E37').Value = '  +2.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01608'
$ws.Range('E38').Value = '  -0.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8660'
$ws.Range('E39').Value = '  +2.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.820'
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('D41').Value = '1.049.73'
$ws.Range('E41').Value = '  +1.87%  '
$ws.Range('E42').Value = '  +0.76%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '104.49'
$ws.Range('E43').Value = '  +1.62%  '
$ws.Range('D44').Value = '1.803.33'
$ws.Range('E44').Value = '  +0.33%  '
$ws.Range('E45').Value = '  +0.99%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₈106'
$ws.Range('E46').Value = '  -5.21%  '
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.009'
$ws.Range('E47').Value = '  +0.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4384'
$ws.Range('E48').Value = '  +1.84%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.035'
$ws.Range('E49').Value = '  +2.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05164'
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.416'
$ws.Range('E51').Value = '  -3.91%  '
